# Updates the "data" template workbook:
#  - removes the embedded instruction/metadata rows (rows 2-5) from the
#    "data" sheet, leaving only the header row
#  - adds a new "lifestage" column description row to the "instructions"
#    sheet (inserted right after "effect", before "simple_lifestage"),
#    and updates the "simple_lifestage" row's description/allowed values
#    to reflect that it is now derived from "lifestage"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "data": drop the old embedded metadata rows (2-5), keep header,
# fix the bogus "name" header to "latin_name", and insert a new
# "lifestage" column (after "effect", before "simple_lifestage").
# ---------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Rows.Item(2).Resize(4).Delete()

$dataSheet.Range("A1").Value = "latin_name"
$dataSheet.Range("B1").Value = "endpoint"
$dataSheet.Range("C1").Value = "effect"
$dataSheet.Range("D1").Value = "lifestage"
$dataSheet.Range("E1").Value = "simple_lifestage"
$dataSheet.Range("F1").Value = "effect_conc_mg.L"
$dataSheet.Range("G1").Value = "duration_hrs"
$dataSheet.Range("H1").Value = "trophic_group"
$dataSheet.Range("I1").Value = "ecological_group"
$dataSheet.Range("J1").Value = "species_present_in_bc"

# ---------------------------------------------------------------
# Sheet "instructions": insert a new row for "lifestage" after the
# "effect" row (row 4), then update the "simple_lifestage" row that
# follows it.
# ---------------------------------------------------------------
$instrSheet = $wb.Worksheets.Item("instructions")

# Insert a new blank row at position 5 (pushes simple_lifestage etc. down)
$instrSheet.Rows.Item(5).Insert()

$instrSheet.Range("A5").Value = "lifestage"
$instrSheet.Range("B5").Value = "The lifestage of the species during the test."
$instrSheet.Range("C5").Value = "Any lifestage"
$instrSheet.Range("D5").Value = "egg"
# Leave E5 empty, but materialize an explicit (empty) cell record so the
# row matches the shape of its siblings (which all carry a trailing
# empty "extra_details" cell).
$instrSheet.Range("E5").ClearContents()
$instrSheet.Range("E5").Style = "Normal"

# simple_lifestage is now row 6 - update its description/allowed_values/example
$instrSheet.Range("A6").Value = "simple_lifestage"
$instrSheet.Range("B6").Value = "A simplified lifestage that is based on the lifestage."
$instrSheet.Range("C6").Value = "One of els, adult, or juvenile"
$instrSheet.Range("D6").Value = "els"
$instrSheet.Range("E6").Value = "early life stage (els)"
